# Apply the "ranking_Aug-2023" update:
#  - Swap the K1/L1 header labels ("h-index (5 years)" <-> "i10-index")
#  - Swap the K/L data values for each data row (2..12)
#  - Renumber column A (rank) sequentially from 1 to 11 for rows 2..12

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap header labels in K1 and L1 ---
$k1 = $ws.Range("K1").Value2
$l1 = $ws.Range("L1").Value2
$ws.Range("K1").Value2 = $l1
$ws.Range("L1").Value2 = $k1

# --- Swap K/L data values and renumber column A for each data row ---
$newRank = 1
for ($row = 2; $row -le 12; $row++) {
    $kCell = $ws.Cells.Item($row, 11)  # column K
    $lCell = $ws.Cells.Item($row, 12)  # column L

    $kVal = $kCell.Value2
    $lVal = $lCell.Value2

    $kCell.Value2 = $lVal
    $lCell.Value2 = $kVal

    $ws.Cells.Item($row, 1).Value2 = $newRank
    $newRank++
}
